$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder "Periodo Mora" values in E16:E21 from descending (2303..2210) to ascending (2210..2303)
$ws.Range("E16").Value = "2210"
$ws.Range("E17").Value = "2211"
$ws.Range("E18").Value = "2212"
$ws.Range("E19").Value = "2301"
$ws.Range("E20").Value = "2302"
$ws.Range("E21").Value = "2303"

# Update "Valor Mora" values in G16:G21 from 1,000,000 to 1,300,000
$ws.Range("G16:G21").Value = 1300000
